$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data edits
$ws.Range("A3").Value = "Admin"
$ws.Range("A4").Value = "admin"
$ws.Range("B4").Value = "pass"

# Unify formatting: B2:B4 previously carried a distinct (font-applied)
# bordered/left-aligned style; make them match column A's bordered
# left-aligned style so the redundant style entry collapses away.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("B2:B4").PasteSpecial(-4122) | Out-Null

# Restore the active selection like the authored workbook
$ws.Range("B15").Select() | Out-Null
